$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.281.79"
$ws.Cells.Item(2, 5).Value = "  +0.15%  "
$ws.Cells.Item(3, 4).Value = "1.599.79"
$ws.Cells.Item(3, 5).Value = "  +0.69%  "
$ws.Cells.Item(4, 5).Value = "  +0.18%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "212.68"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.39%  "
$ws.Cells.Item(6, 5).Value = "  -0.18%  "
$ws.Cells.Item(7, 5).Value = "  +0.22%  "
$ws.Cells.Item(8, 5).Value = "  -0.51%  "
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.0606"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +0.06%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "18.97"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -1.47%  "
$ws.Cells.Item(11, 5).Value = "  +1.02%  "
$ws.Cells.Item(12, 5).Value = "  +0.82%  "
$ws.Cells.Item(13, 4).Value = "1.604.89"
$ws.Cells.Item(13, 5).Value = "  +1.52%  "
$ws.Cells.Item(14, 5).Value = "  -0.24%  "
$ws.Cells.Item(15, 5).Value = "  -2.19%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "63.65"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.70%  "
$ws.Cells.Item(17, 4).Value = "26.284.85"
$ws.Cells.Item(17, 5).Value = "  +0.23%  "
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "229.52"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +6.61%  "
$ws.Cells.Item(19, 5).Value = "  -0.64%  "
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "7.60"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +3.49%  "
$ws.Cells.Item(21, 5).Value = "  +0.06%  "
$ws.Cells.Item(22, 5).Value = "  +0.81%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "2.17"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.18%  "
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "145.55"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.04%  "
$ws.Cells.Item(26, 5).Value = "  +0.20%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "6.95"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.61%  "
$ws.Cells.Item(28, 5).Value = "  +1.01%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "15.42"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +1.74%  "
$ws.Cells.Item(30, 5).Value = "  -0.80%  "
$ws.Cells.Item(31, 5).Value = "  +0.79%  "
$ws.Cells.Item(32, 5).Value = "  -0.07%  "
$ws.Cells.Item(33, 4).Value = "1.445.26"
$ws.Cells.Item(33, 5).Value = "  +5.93%  "
$ws.Cells.Item(34, 5).Value = "  +0.89%  "
$ws.Cells.Item(35, 5).Value = "  -0.35%  "
$ws.Cells.Item(36, 5).Value = "  +0.59%  "
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "0.569"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -2.28%  "
$ws.Cells.Item(38, 5).Value = "  -1.24%  "
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.822"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +0.25%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "5.78"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -0.89%  "
$ws.Cells.Item(41, 5).Value = "  +0.24%  "
$ws.Cells.Item(42, 5).Value = "  +2.16%  "
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "0.923"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.02%  "
$ws.Cells.Item(44, 4).Value = "1.738.83"
$ws.Cells.Item(44, 5).Value = "  +0.88%  "
$ws.Cells.Item(45, 5).Value = "  -1.27%  "
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "60.70"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -0.53%  "
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "87.34"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +1.30%  "
$ws.Cells.Item(48, 5).Value = "  +0.11%  "
$ws.Cells.Item(49, 5).Value = "  +0.00%  "
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "0.0951"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -3.38%  "
$ws.Cells.Item(51, 5).Value = "  +0.21%  "
